# "Add files via upload" - append three new employee records
# (matricula / nome / uuid) to the bottom of the worksheet, using the
# same small Arial font already used for the preceding rows, and move
# the active-cell selection the way Excel would after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New records to append (matricula, nome, uuid)
$novos = @(
    @{ A = 5289; B = "ROSY KELLY CORDEIRO"; C = "a5bc41fc-7684-44c5-9817-a6b1c315c403" },
    @{ A = 5286; B = "ANA LUCIA MOTA LINO";  C = "520b82a5-312f-4441-8434-0b2af059dc91" },
    @{ A = 5287; B = "ROSALINDA SOARES";     C = "205f5ec9-a3b2-4deb-89fc-dba13e1288f5" }
)

$startRow = 160
for ($i = 0; $i -lt $novos.Count; $i++) {
    $row = $startRow + $i
    $rec = $novos[$i]

    $ws.Cells.Item($row, 1).Value = $rec.A
    $ws.Cells.Item($row, 2).Value = $rec.B
    $ws.Cells.Item($row, 3).Value = $rec.C
}

# Format the first new name cell with the small Arial font ...
$firstName = $ws.Cells.Item($startRow, 2)
$firstName.Font.Name = "Arial"
$firstName.Font.Size = 6
$firstName.Font.Color = 1907741

# ... then copy that formatting onto the rest of the new name cells so
# every added row shares a single, identical cell style.
if ($novos.Count -gt 1) {
    $firstName.Copy()
    $restNames = $ws.Range($ws.Cells.Item($startRow + 1, 2), $ws.Cells.Item($startRow + $novos.Count - 1, 2))
    $restNames.PasteSpecial(-4122)
}

# Move the selection the way Excel would after adding the rows
$ws.Range("B163").Select()
